$d = $word.ActiveDocument

# The target text " business analytics. " lives in a single paragraph and is
# split (in the original) across two runs:
#   run A: " business "   -> should become " business"
#   run B: "analytics. "  -> should become ". "
# i.e. the phrase " analytics" (the space before it plus the word itself) is
# removed, turning "Enhancing business analytics. " into
# "Enhancing business. ".

$locate = $d.Content
$found = $locate.Find.Execute(" business analytics. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate ' business analytics. ' in the document."
}

$matchStart = $locate.Start
$matchEnd = $locate.End

# Delete " analytics" (the space that follows "business" together with the
# word "analytics") leaving " business. " in its place.
$toDelete = $d.Range($matchStart + 9, $matchStart + 19)
$toDelete.Delete()

# The deletion above merges the surrounding runs into a single run (this
# runtime always re-flows a paragraph's runs when text is edited). Re-impose
# a run boundary exactly where the original two runs used to meet by nudging
# (and restoring) a formatting property on each side -- this forces the
# engine to split the run again without altering any visible formatting.
$newPeriodRun = $d.Range($matchStart + 9, $matchStart + 11)   # ". "
$newPeriodRun.Font.Bold = 1
$newPeriodRun.Font.Bold = 0

$newBusinessRun = $d.Range($matchStart, $matchStart + 9)      # " business"
$newBusinessRun.Font.Bold = 1
$newBusinessRun.Font.Bold = 0
